# "feat: add 2022-Q3 data"
#
# Before:  总计 | 2022-Q1 | 2021-Q4
# After:   总计 | 2022-Q3 | 2022-Q1 | 2021-Q4
#
#  - 总计 (summary) sheet gains a new row for 2022-Q3, pushing the
#    2022-Q1 / 2021-Q4 rows down by one.
#  - A new "2022-Q3" fund-holdings sheet appears with fresh figures.
#  - The old "2022-Q1" fund-holdings sheet (and its figures) is kept,
#    just shifted one position to the right.
#  - "2021-Q4" is left completely untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 sheet: insert a 2022-Q3 row above the existing 2022-Q1 row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(3, 1).Copy($summary.Cells.Item(2, 1))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 0.13
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2

# ---------------------------------------------------------------------
# 2) Duplicate the current "2022-Q1" sheet (so its original data/styling
#    is preserved verbatim in the new, shifted-right "2022-Q1" sheet),
#    then turn the original sheet into "2022-Q3" and refresh its figures.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1Index = $q1.Index

$q1.Copy($null, $q1)               # places an exact duplicate right after $q1

$q1.Name = "2022-Q3"                # the original sheet becomes 2022-Q3
$newQ1 = $wb.Worksheets.Item($q1Index + 1)
$newQ1.Name = "2022-Q1"             # the duplicate keeps the old data, renamed back

# New 2022-Q3 fund-holding figures (row layout: code, name, scale, total
# stock position, position share, holding value, position rank).
$q3Data = @(
    @("'000369", "广发全球医疗保健（QDII）人民币A", "'2.76", "'83.19", "'2.42", "'0.0668", 10),
    @("'000370", "广发全球医疗保健（QDII）美元A", "'2.75", "'83.19", "'2.42", "'0.0666", 10),
    @("'016280", "广发全球医疗保健（QDII）人民币C", "'0.02", "'83.19", "'2.42", "'0.0005", 10),
    @("'016281", "广发全球医疗保健（QDII）美元C", "'0.02", "'83.19", "'2.42", "'0.0005", 10)
)

for ($i = 0; $i -lt $q3Data.Length; $i++) {
    $row = $i + 2
    $rowData = $q3Data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 2
        $q1.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

Write-Host "2022-Q3 data added"
